$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.678.19"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.962.29"
$ws.Range("E3").Value = "  +1.62%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'244.40"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("D7").Value = "'60.50"
$ws.Range("E7").Value = "  +8.47%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.376"
$ws.Range("E9").Value = "  +5.75%  "
$ws.Range("D10").Value = "'0.0793"
$ws.Range("E10").Value = "  -4.63%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E12").Value = "  +7.49%  "
$ws.Range("D13").Value = "'0.844"
$ws.Range("E13").Value = "  +6.12%  "
$ws.Range("D14").Value = "2.249.49"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "'21.66"
$ws.Range("E15").Value = "  +4.90%  "
$ws.Range("E16").Value = "  +4.40%  "
$ws.Range("D17").Value = "1.957.65"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "36.611.68"
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "'69.84"
$ws.Range("E19").Value = "  +1.89%  "
$ws.Range("D20").Value = "0.0₃0853"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("D21").Value = "'230.26"
$ws.Range("E21").Value = "  +1.98%  "
$ws.Range("D22").Value = "'5.10"
$ws.Range("E22").Value = "  +3.88%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +6.95%  "
$ws.Range("E25").Value = "  +4.91%  "
$ws.Range("E26").Value = "  +12.57%  "
$ws.Range("D27").Value = "'9.25"
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("D28").Value = "'160.82"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("E29").Value = "  +2.60%  "
$ws.Range("E30").Value = "  +11.00%  "
$ws.Range("E31").Value = "  +3.12%  "
$ws.Range("E32").Value = "  +5.84%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("D34").Value = "'4.43"
$ws.Range("E34").Value = "  +7.90%  "
$ws.Range("E35").Value = "  +21.79%  "
$ws.Range("E36").Value = "  +8.41%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").Value = "'5.64"
$ws.Range("E39").Value = "  -5.92%  "
$ws.Range("E40").Value = "  +2.18%  "
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("E42").Value = "  +3.51%  "
$ws.Range("E43").Value = "  +1.56%  "
$ws.Range("D44").Value = "'16.14"
$ws.Range("E44").Value = "  +4.83%  "
$ws.Range("D45").Value = "1.371.64"
$ws.Range("E45").Value = "  +3.54%  "
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("D47").Value = "'88.46"
$ws.Range("E47").Value = "  +4.80%  "
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").Value = "'44.54"
$ws.Range("E50").Value = "  +3.39%  "
$ws.Range("D51").Value = "2.139.73"
$ws.Range("E51").Value = "  +1.61%  "
